# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (fund holdings for that quarter) right
# after the "总计" (totals) sheet, and updates the "总计" sheet so its
# summary table includes the new quarter as its first data row.
#
# All of the other quarter sheets ("2022-Q3", "2022-Q1", "2021-Q3",
# "2021-Q1") keep their own data untouched - they simply shift one tab to
# the right to make room for the new sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# Header row (bold, centered, thin box border - matches the other fund
# sheets).
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    $q4.Cells.Item(1, $c).Value = $headers[$c - 2]
}

# Data rows. Columns B (fund code) and D:G (scale / position figures) are
# kept as text, just like the existing sheets, so codes such as "002295"
# keep their leading zero and the figures aren't reformatted as numbers.
$q4Rows = @(
    @(0, "002295", "广发稳安灵活配置混合A", "1.51", "80.18", "4.31", "0.0651", 9),
    @(1, "008604", "广发稳安灵活配置混合C", "0.01", "80.18", "4.31", "0.0004", 9)
)

$rowIdx = 2
foreach ($row in $q4Rows) {
    $aCell = $q4.Cells.Item($rowIdx, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
    $aCell.Value = $row[0]

    $bCell = $q4.Cells.Item($rowIdx, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $q4.Cells.Item($rowIdx, 3).Value = $row[2]

    for ($c = 4; $c -le 7; $c++) {
        $cell = $q4.Cells.Item($rowIdx, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
    }

    $q4.Cells.Item($rowIdx, 8).Value = $row[7]

    $rowIdx = $rowIdx + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add the 2022-Q4 row at the top of
#    the data and shift the older rows down, finishing with the 2021-Q1
#    row that used to be last.
# ---------------------------------------------------------------------
$totalRows = @(
    @("2022-Q4", 2, 0.07),
    @("2022-Q3", 2, 0.07),
    @("2022-Q1", 1, 0.08),
    @("2021-Q3", 1, 0.09),
    @("2021-Q1", 4, 0.1)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# New row 6 needs the same index-column styling (bold + border) as the
# existing A2:A5 cells.
$a6 = $totalSheet.Cells.Item(6, 1)
$a6.Font.Bold = $true
$a6.HorizontalAlignment = -4108
$a6.VerticalAlignment = -4160
$a6.Borders.LineStyle = 1
$a6.Value = 4

# ---------------------------------------------------------------------
# 3. Keep "2021-Q1" as the selected/active sheet, matching the workbook
#    before this edit.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count())
$lastSheet.Activate()
